$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for all data rows (2-98):
# every value moves from 45205 (2023-10-06) to 45206 (2023-10-07).
for ($row = 2; $row -le 98; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45205) {
        $cell.Value = 45206
    }
}
